$d = $word.ActiveDocument

# 1) Expand "as well as a csv file ... final order." into the new two-sentence
#    wording, inserting the "turn by turn directions" sentence and the
#    "has also been attached" tail.
$find1 = "as well as a csv file of the coordinates in their final order."
$repl1 = "as well as turn by turn directions. Additionally, a csv file of the coordinates in their final order has also been attached."
$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2) | Out-Null

# 2) Collapse the double space that used to separate "minutes." from "Leg 3"
#    down to a single space.
$find2 = "55 minutes.  Leg 3"
$repl2 = "55 minutes. Leg 3"
$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2) | Out-Null

# 3) Move the (hidden) "_GoBack" bookmark that Word maintains for "last edit
#    location" so that it sits right after "...has also been attached" and
#    before the following period, matching where the real edit happened.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$rng = $d.Content
$rng.Find.Execute("has also been attached") | Out-Null
$markRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $markRange) | Out-Null
